$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the empty inline-string placeholder cells in row 2 (A2-F2, I2-N2)
# leaving only G2 (Date) and H2 (Time) populated, as in the original.
$ws.Range("A2:F2").ClearContents()
$ws.Range("I2:N2").ClearContents()

# The Freq./Readability/Signal-strength columns hold numeric-looking log
# entries ("14.250", "5", "9") that must stay literal text, not get
# auto-converted to numbers - so mark them as Text before typing them in.
$ws.Range("I3:I4").NumberFormat = "@"
$ws.Range("K3:L4").NumberFormat = "@"

# Row 3 - first new contact entry
$row3 = @("ON4CJK", "Jose", "JO11ub", "ON1DDR", "Andre", "JO11ub", "18-10-2024", "10:46", "14.250", "SSB", "5", "9", "/", "via RTL-SDR")
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

# Row 4 - second new contact entry
$row4 = @("ON4CJK", "jose", "JO11ub", "ON1DDR", "Jose", "JO11ds", "18-10-2024", "11:19", "14.250", "SSB", "5", "9", "/", "Mary Islands")
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}
